$d = $word.ActiveDocument

# --- Edit 1: Title paragraph text "Dummy " -> "Dummy 2 " -------------------
# (insert " 2" immediately after "Dummy", before the existing trailing
# space / "Trial " run, so the title reads "Dummy 2 Trial 2 Test Case ...")
$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Dummy", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$titleRange.Collapse(0)
$titleRange.InsertAfter(" 2")

# --- Edit 2: Append a new paragraph "MODIFIED BY BACKEND" at the very end --
# (right after the "Expected Result: Accepted - Valid password" paragraph,
# before the section properties)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "MODIFIED BY BACKEND"
